{"js": "// Word JS API (Office.js) edit script.\n//\n// Source change (from the commit's XML diff):\n//  1. Remove the sentence \"Have 2 wards that don't have data are Thu Thiem\n//     (27118) and An Loi Dong (27115) because of the city government\n//     re-planning residential areas.\" \u2014 the run is deleted but the\n//     paragraph that carried it is left in place (now empty).\n//  2. Right after the paragraph that follows the INCLUDEPICTURE\n//     image (an already-empty paragraph), insert a brand-new paragraph\n//     containing: \"Note: population in Thu Thiem and An Loi Dong: 2021\n//     using population in 2019, 2030: using population in population\n//     planning of this ward.\"\n\nconst body = context.document.body;\n\n// ---------------------------------------------------------------------\n// Change 1 \u2014 delete the \"Have 2 wards...\" sentence, keep its paragraph.\n// ---------------------------------------------------------------------\nconst oldSentence =\n  \"Have 2 wards that don\\u2019t have data are Thu Thiem (27118) and An Loi \" +\n  \"Dong (27115) because of the city government re-planning residential areas.\";\n\nconst matches = body.search(oldSentence, { matchCase: true });\nmatches.load(\"items\");\nawait context.sync();\n\nif (matches.items.length > 0) {\n  // Deleting the matched Range (rather than the whole paragraph's range)\n  // removes just the run and leaves the enclosing <w:p> intact \u2014 matching\n  // the diff, which keeps the now-empty paragraph.\n  matches.items[0].delete();\n  await context.sync();\n}\n\n// ---------------------------------------------------------------------\n// Change 2 \u2014 add the new \"Note: ...\" paragraph after the paragraph that\n// follows the inline picture (and before the document's trailing empty\n// paragraph).\n// ---------------------------------------------------------------------\nconst noteSentence =\n  \"Note: population in Thu Thiem and An Loi Dong: 2021 using population in \" +\n  \"2019, 2030: using population in population planning of this ward.\";\n\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Load each paragraph's inline pictures so we can find the one holding the\n// INCLUDEPICTURE image.\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  paragraphs.items[i].inlinePictures.load(\"items\");\n}\nawait context.sync();\n\nlet pictureParagraphIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].inlinePictures.items.length > 0) {\n    pictureParagraphIndex = i;\n    break;\n  }\n}\n\n// The paragraph right after the picture paragraph is the (empty) one the\n// new \"Note:\" paragraph must follow.\nconst anchorParagraph =\n  pictureParagraphIndex >= 0 && pictureParagraphIndex + 1 < paragraphs.items.length\n    ? paragraphs.items[pictureParagraphIndex + 1]\n    : paragraphs.items[paragraphs.items.length - 2]; // fallback: 2nd-to-last paragraph\n\nanchorParagraph.insertParagraph(noteSentence, \"After\");\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n#\n# Source change (from the commit's XML diff):\n#  1. Remove the sentence \"Have 2 wards that don't have data are Thu Thiem\n#     (27118) and An Loi Dong (27115) because of the city government\n#     re-planning residential areas.\" \u2014 only the run is removed; the\n#     paragraph that carried it is left in place (now empty).\n#  2. Right after the paragraph that follows the INCLUDEPICTURE image (an\n#     already-empty paragraph), insert a brand-new paragraph containing:\n#     \"Note: population in Thu Thiem and An Loi Dong: 2021 using\n#     population in 2019, 2030: using population in population planning\n#     of this ward.\"\n\n$d = $word.ActiveDocument\n\n# ---------------------------------------------------------------------\n# Change 1 - delete the \"Have 2 wards...\" sentence, keep its paragraph.\n# ---------------------------------------------------------------------\n$oldSentence = \"Have 2 wards that don\u2019t have data are Thu Thiem (27118) and An Loi Dong (27115) because of the city government re-planning residential areas.\"\n\n$findRange = $d.Content\n$found = $findRange.Find.Execute($oldSentence)\n\nif ($found) {\n    # Grow the hit to its whole enclosing paragraph, then back off one\n    # character so the trailing paragraph mark is excluded \u2014 deleting that\n    # trimmed range removes just the run and keeps the (now empty) <w:p>,\n    # matching the diff.\n    $findRange.Expand(4) | Out-Null    # wdParagraph = 4\n    $findRange.MoveEnd(1, -1) | Out-Null   # wdCharacter = 1\n    $findRange.Delete() | Out-Null\n}\n\n# ---------------------------------------------------------------------\n# Change 2 - add the new \"Note: ...\" paragraph after the paragraph that\n# follows the inline picture (and before the document's trailing empty\n# paragraph).\n# ---------------------------------------------------------------------\n$noteSentence = \"Note: population in Thu Thiem and An Loi Dong: 2021 using population in 2019, 2030: using population in population planning of this ward.\"\n\n$paraCount = $d.Paragraphs.Count\n$pictureParaIndex = -1\nfor ($i = 1; $i -le $paraCount; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.InlineShapes.Count -gt 0) {\n        $pictureParaIndex = $i\n        break\n    }\n}\n\nif ($pictureParaIndex -gt 0 -and ($pictureParaIndex + 1) -le $d.Paragraphs.Count) {\n    $anchorParaIndex = $pictureParaIndex + 1\n} else {\n    # Fallback: the paragraph just before the document's final paragraph.\n    $anchorParaIndex = $d.Paragraphs.Count - 1\n}\n\n$anchorPara = $d.Paragraphs.Item($anchorParaIndex)\n$anchorPara.Range.InsertParagraphAfter() | Out-Null\n$d.Paragraphs.Item($anchorParaIndex + 1).Range.Text = $noteSentence\n"}
